$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.741.35'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '2.341.49'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.17'
$ws.Range("E5").Value = '  +5.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.95'
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.538'
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("D9").Value = '2.362.38'
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("E10").Value = '  +1.94%  '
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.40'
$ws.Range("E12").Value = '  +1.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +4.38%  '
$ws.Range("D14").Value = '2.782.88'
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.56'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = '57.807.79'
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("D18").Value = '2.382.98'
$ws.Range("E18").Value = '  +2.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.58'
$ws.Range("E19").Value = '  +1.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '335.26'
$ws.Range("E20").Value = '  +3.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.21'
$ws.Range("E21").Value = '  +1.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.76'
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.13'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("E25").Value = '  +3.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.45'
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.43'
$ws.Range("E28").Value = '  +8.78%  '
$ws.Range("E29").Value = '  +4.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.63'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("D31").Value = '0.0₃0738'
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.18'
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.59'
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.03'
$ws.Range("E34").Value = '  +15.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.994'
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.15'
$ws.Range("E38").Value = '  +5.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.62'
$ws.Range("E39").Value = '  +3.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.43'
$ws.Range("E40").Value = '  +2.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '150.21'
$ws.Range("E41").Value = '  -3.00%  '
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("E43").Value = '  +1.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '284.20'
$ws.Range("E44").Value = '  +1.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.34'
$ws.Range("E45").Value = '  +6.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0932'
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0506'
$ws.Range("E47").Value = '  +1.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.563'
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0219'
$ws.Range("E49").Value = '  +1.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.62'
$ws.Range("E50").Value = '  +2.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.381'
$ws.Range("E51").Value = '  -0.23%  '
